$wb = $excel.ActiveWorkbook

# --- VoiceContinuity sheet (7th sheet): replace the per-flavour default rows
# with a single "configurationByDefault" row driven by a CSV of values, and
# clear the now-unused rows 3-8. ---
$wsVoiceContinuity = $wb.Worksheets.Item(7)

$wsVoiceContinuity.Range("A2").Value = "configurationByDefault"
$wsVoiceContinuity.Range("B2").Value = "New,Existing,Phone Line,475153060,NotApplicable,NotApplicable,Voice Continuity 1"

# B2 now picks up the plain "text" column style (matching column B's default
# style) instead of the old row style.
$wsVoiceContinuity.Range("B2").Style = "Normal"
$wsVoiceContinuity.Range("B2").NumberFormat = "@"

$wsVoiceContinuity.Range("A3:B8").ClearContents()

# --- Update cursor/selection position on a few sheets (view-only changes) ---
$wb.Worksheets.Item(1).Range("B12").Select()
$wb.Worksheets.Item(2).Range("B3").Select()
$wb.Worksheets.Item(3).Range("B39").Select()

# Activate VoiceContinuity last so it becomes the saved active tab, matching
# the new tabSelected/activeTab state, and set its new selection.
$wsVoiceContinuity.Range("A3").Select()
